# Generate Report for Handback
#
# For both the "zh-cn" and "de-de" localization status sheets:
#  - Mark the Status column (B) as "Handed back: in sync with en-US"
#  - Populate "Latest Target File" (E) and "Latest Handback File" (F) with
#    hyperlinks that mirror the existing "Source File Name" (A) and
#    "Latest Handoff File" (C) hyperlinks
#  - Stamp "Latest Handback DateTime" (G) with the actual handback time

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# The Overview sheet mirrors the same "Status" text for each locale/file
# combination, so it needs to be brought in sync as well.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

function Update-HandbackSheet {
    param(
        [string]$SheetName,
        [string]$MdDisplay,
        [string]$MdAddress,
        [string]$XlfDisplay,
        [string]$XlfAddress,
        [string]$MdDisplay2,
        [string]$MdAddress2,
        [string]$XlfDisplay2,
        [string]$XlfAddress2,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Row 2 - c64627a0-f454-4aeb-ad95-7f17527d6634
    $ws.Range("B2").Value = $statusHandedBack
    $ws.Hyperlinks.Add($ws.Range("E2"), $MdAddress, [Type]::Missing, [Type]::Missing, $MdDisplay)
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfAddress, [Type]::Missing, [Type]::Missing, $XlfDisplay)
    $ws.Range("G2").Value = $HandbackDateTime

    # Row 3 - 0926ad3d-e71c-48d7-851d-d96698f3d4df
    $ws.Range("B3").Value = $statusHandedBack
    $ws.Hyperlinks.Add($ws.Range("E3"), $MdAddress2, [Type]::Missing, [Type]::Missing, $MdDisplay2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $XlfAddress2, [Type]::Missing, [Type]::Missing, $XlfDisplay2)
    $ws.Range("G3").Value = $HandbackDateTime
}

Update-HandbackSheet `
    "zh-cn" `
    "c64627a0-f454-4aeb-ad95-7f17527d6634.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/07f6abcd9e5ec1fa07255a5299b7a0c1931e145a/e2e/c64627a0-f454-4aeb-ad95-7f17527d6634.md" `
    "c64627a0-f454-4aeb-ad95-7f17527d6634.6c9b2070688f53a89a700c399e9ad5312c8bf40f.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1feec7d24f55a0f0dede78feaba062b8a741c40b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c64627a0-f454-4aeb-ad95-7f17527d6634.6c9b2070688f53a89a700c399e9ad5312c8bf40f.zh-cn.xlf" `
    "0926ad3d-e71c-48d7-851d-d96698f3d4df.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/07f6abcd9e5ec1fa07255a5299b7a0c1931e145a/e2e/0926ad3d-e71c-48d7-851d-d96698f3d4df.md" `
    "0926ad3d-e71c-48d7-851d-d96698f3d4df.7713a32d8ba09a0caebfa5f88b2c869af20a614e.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1feec7d24f55a0f0dede78feaba062b8a741c40b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0926ad3d-e71c-48d7-851d-d96698f3d4df.7713a32d8ba09a0caebfa5f88b2c869af20a614e.zh-cn.xlf" `
    "2016-03-03 10:58:31"

Update-HandbackSheet `
    "de-de" `
    "c64627a0-f454-4aeb-ad95-7f17527d6634.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/07f6abcd9e5ec1fa07255a5299b7a0c1931e145a/e2e/c64627a0-f454-4aeb-ad95-7f17527d6634.md" `
    "c64627a0-f454-4aeb-ad95-7f17527d6634.6c9b2070688f53a89a700c399e9ad5312c8bf40f.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78a51e7d293737e18d1fd08ee71106c88f6fe4b3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c64627a0-f454-4aeb-ad95-7f17527d6634.6c9b2070688f53a89a700c399e9ad5312c8bf40f.de-de.xlf" `
    "0926ad3d-e71c-48d7-851d-d96698f3d4df.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/07f6abcd9e5ec1fa07255a5299b7a0c1931e145a/e2e/0926ad3d-e71c-48d7-851d-d96698f3d4df.md" `
    "0926ad3d-e71c-48d7-851d-d96698f3d4df.7713a32d8ba09a0caebfa5f88b2c869af20a614e.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/78a51e7d293737e18d1fd08ee71106c88f6fe4b3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0926ad3d-e71c-48d7-851d-d96698f3d4df.7713a32d8ba09a0caebfa5f88b2c869af20a614e.de-de.xlf" `
    "2016-03-03 10:58:55"
